$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Year / Month / Data Availability / Data Missing
$ws.Range("A2").Value = 2024
$ws.Range("B2").Value = "DEC"
$ws.Range("C2").Value = "31/12-01/12"
$ws.Range("D2").Value = "-"

# Station Name stays the same
$ws.Range("E2").Value = "Akurana"

# Columns F through AG become checkmarks, except M, Y, Z which become "-"
$checkCols = @("F","G","H","I","J","K","L","N","O","P","Q","R","S","T","U","V","W","X","AA","AB","AC","AD","AE","AF","AG")
foreach ($col in $checkCols) {
    $ws.Range("$col" + "2").Value = "✓"
}

$dashCols = @("M","Y","Z","AH","AI","AJ")
foreach ($col in $dashCols) {
    $ws.Range("$col" + "2").Value = "-"
}
